$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 6034
$ws.Range("K3").Value = 6218
$ws.Range("F4").Value = 1579
$ws.Range("J4").Value = 1498
$ws.Range("K4").Value = 1301
$ws.Range("K5").Value = 440
$ws.Range("K6").Value = 6833
$ws.Range("F7").Value = 19401
$ws.Range("J7").Value = 23204
$ws.Range("K7").Value = 20826

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K2").Value = 183
$ws.Range("K7").Value = 608
$ws.Range("K8").Value = 1369
$ws.Range("K9").Value = 91
$ws.Range("K11").Value = 390
$ws.Range("K18").Value = 138
$ws.Range("K19").Value = 603
$ws.Range("K20").Value = 495
$ws.Range("K23").Value = 212
$ws.Range("K24").Value = 61
$ws.Range("K25").Value = 98
$ws.Range("K27").Value = 194
$ws.Range("K33").Value = 903
$ws.Range("K34").Value = 120
$ws.Range("K37").Value = 710
$ws.Range("K38").Value = 20
$ws.Range("K39").Value = 26
$ws.Range("K41").Value = 146
$ws.Range("K42").Value = 774
$ws.Range("K43").Value = 175
$ws.Range("K47").Value = 145
$ws.Range("K48").Value = 262
$ws.Range("K51").Value = 268
$ws.Range("K52").Value = 544
$ws.Range("K53").Value = 267
$ws.Range("K55").Value = 229
$ws.Range("F63").Value = 158
$ws.Range("J63").Value = 101
$ws.Range("K63").Value = 61
$ws.Range("K65").Value = 490
$ws.Range("K67").Value = 811
$ws.Range("K72").Value = 102
$ws.Range("K73").Value = 184
$ws.Range("K76").Value = 281
$ws.Range("K77").Value = 145
$ws.Range("K82").Value = 22
$ws.Range("K83").Value = 462
$ws.Range("K85").Value = 966
$ws.Range("K88").Value = 224
$ws.Range("K89").Value = 303
$ws.Range("K90").Value = 192
$ws.Range("K91").Value = 236
$ws.Range("K93").Value = 77
$ws.Range("K94").Value = 280
$ws.Range("K95").Value = 350
$ws.Range("K96").Value = 219
$ws.Range("F101").Value = 19401
$ws.Range("J101").Value = 23204
$ws.Range("K101").Value = 20826

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 68
$ws.Range("K7").Value = 219

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 200
$ws.Range("K7").Value = 608

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range("K2").Value = 138
$ws.Range("K3").Value = 101
$ws.Range("K6").Value = 125
$ws.Range("K7").Value = 390

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range("K3").Value = 92
$ws.Range("K7").Value = 303

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K3").Value = 335
$ws.Range("K7").Value = 966

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K2").Value = 147
$ws.Range("K4").Value = 30
$ws.Range("K7").Value = 544

$ws = $wb.Worksheets.Item('Logan Square')
$ws.Range("K2").Value = 69
$ws.Range("K7").Value = 267

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 375
$ws.Range("K3").Value = 418
$ws.Range("K6").Value = 463
$ws.Range("K7").Value = 1369

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 159
$ws.Range("K3").Value = 165
$ws.Range("K7").Value = 462

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 330
$ws.Range("K4").Value = 45
$ws.Range("K7").Value = 903

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 124
$ws.Range("K7").Value = 350

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K3").Value = 235
$ws.Range("K5").Value = 31
$ws.Range("K6").Value = 207
$ws.Range("K7").Value = 710

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 162
$ws.Range("K7").Value = 490

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 225
$ws.Range("K3").Value = 291
$ws.Range("K6").Value = 231
$ws.Range("K7").Value = 811

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K4").Value = 37
$ws.Range("K7").Value = 262

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K4").Value = 29
$ws.Range("K6").Value = 194
$ws.Range("K7").Value = 603

$ws = $wb.Worksheets.Item('River North')
$ws.Range("K2").Value = 62
$ws.Range("K3").Value = 54
$ws.Range("K7").Value = 281

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K3").Value = 29
$ws.Range("K7").Value = 146

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 233
$ws.Range("K6").Value = 287
$ws.Range("K7").Value = 774

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range("K2").Value = 72
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K2").Value = 23
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range("K3").Value = 74
$ws.Range("K7").Value = 212

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 60
$ws.Range("K3").Value = 112
$ws.Range("K7").Value = 236

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 164
$ws.Range("K3").Value = 161
$ws.Range("K7").Value = 495

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K3").Value = 45
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item('West Lawn')
$ws.Range("K3").Value = 18
$ws.Range("K7").Value = 77

$ws = $wb.Worksheets.Item('Garfield Ridge')
$ws.Range("K2").Value = 47
$ws.Range("K7").Value = 120

$ws = $wb.Worksheets.Item('West Loop')
$ws.Range("K3").Value = 56
$ws.Range("K7").Value = 280

$ws = $wb.Worksheets.Item('East Side')
$ws.Range("K2").Value = 35
$ws.Range("K5").Value = 3
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K6").Value = 47
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("K5").Value = 16
$ws.Range("K6").Value = 26

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K6").Value = 23
$ws.Range("K7").Value = 91

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K3").Value = 47
$ws.Range("K7").Value = 184

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range("K3").Value = 47
$ws.Range("K7").Value = 183

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K2").Value = 57
$ws.Range("K6").Value = 93
$ws.Range("K7").Value = 224

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K2").Value = 53
$ws.Range("K4").Value = 24
$ws.Range("K7").Value = 194

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Range("K2").Value = 71
$ws.Range("K7").Value = 192

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K2").Value = 76
$ws.Range("K7").Value = 268

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range("K2").Value = 35
$ws.Range("K7").Value = 175

$ws = $wb.Worksheets.Item('Old Town')
$ws.Range("K6").Value = 52
$ws.Range("K7").Value = 102

$ws = $wb.Worksheets.Item('Sheffield & DePaul')
$ws.Range("K2").Value = 4
$ws.Range("K6").Value = 22

$ws = $wb.Worksheets.Item('Riverdale')
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 145

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Range("K5").Value = 10
$ws.Range("K6").Value = 20

Write-Host "Applied 164 cell updates across 45 sheets"
